$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set values for "model export" (D) and "model picture" (E) columns for
# newly-added ships, plus the "no model" note (F) for a few rows.

$ws.Cells.Item(11, 5).Value = 1   # E11
$ws.Cells.Item(12, 5).Value = 1   # E12

$ws.Cells.Item(14, 4).Value = 1   # D14
$ws.Cells.Item(14, 5).Value = 1   # E14

$ws.Cells.Item(15, 4).Value = 1   # D15
$ws.Cells.Item(15, 5).Value = 1   # E15

$ws.Cells.Item(16, 4).Value = 1   # D16
$ws.Cells.Item(16, 5).Value = 1   # E16

$ws.Cells.Item(17, 4).Value = 1   # D17
$ws.Cells.Item(17, 5).Value = 1   # E17

$ws.Cells.Item(18, 4).Value = 1   # D18

$ws.Cells.Item(19, 6).Value = "无模型"   # F19

$ws.Cells.Item(20, 4).Value = 1   # D20

$ws.Cells.Item(21, 6).Value = "无模型"   # F21

$ws.Cells.Item(22, 6).Value = "无模型"   # F22

# Update the view state: move the active selection to C28 (this also
# clears the previous topLeftCell scroll offset that pinned the view at
# A4, since the view naturally re-centers on the new selection).
$ws.Range("C28").Select()
